$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8394890427589417
$ws.Range("B1").Value = 3.351017475128174
$ws.Range("C1").Value = 2.451785802841187
$ws.Range("D1").Value = 2.151190996170044
$ws.Range("E1").Value = 1.84245240688324
